# Refresh the cryptos price/volume snapshot (GitHub Actions daily update).
# Price cells that look like plain numbers are prefixed with a leading
# apostrophe so Excel keeps storing them as text (matching the sheet's
# existing inline-string cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.499.68"
$ws.Range("E2").Value = "  -1.81%  "

$ws.Range("D3").Value = "2.427.43"
$ws.Range("E3").Value = "  -2.40%  "

$ws.Range("D5").Value = "'510.36"
$ws.Range("E5").Value = "  -3.08%  "

$ws.Range("D6").Value = "'128.78"
$ws.Range("E6").Value = "  -3.92%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  -2.37%  "

$ws.Range("D9").Value = "2.437.24"
$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("D11").Value = "'0.0945"
$ws.Range("E11").Value = "  -5.71%  "

$ws.Range("D12").Value = "'5.15"
$ws.Range("E12").Value = "  -5.04%  "

$ws.Range("D13").Value = "'0.329"
$ws.Range("E13").Value = "  -4.40%  "

$ws.Range("D14").Value = "2.861.19"
$ws.Range("E14").Value = "  -2.28%  "

$ws.Range("D15").Value = "57.431.51"
$ws.Range("E15").Value = "  -1.75%  "

$ws.Range("D16").Value = "'21.69"
$ws.Range("E16").Value = "  -3.57%  "

$ws.Range("E17").Value = "  -3.81%  "

$ws.Range("D18").Value = "2.435.41"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "'10.39"
$ws.Range("E19").Value = "  -5.20%  "

$ws.Range("D20").Value = "'314.48"
$ws.Range("E20").Value = "  -2.37%  "

$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "'5.62"
$ws.Range("E23").Value = "  -3.68%  "

$ws.Range("D24").Value = "'63.27"
$ws.Range("E24").Value = "  -1.79%  "

$ws.Range("D25").Value = "'0.403"
$ws.Range("E25").Value = "  -2.61%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").Value = "'7.20"
$ws.Range("E28").Value = "  -3.85%  "

$ws.Range("D29").Value = "'170.01"
$ws.Range("E29").Value = "  +2.52%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'6.21"
$ws.Range("E30").Value = "  -3.62%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0715"
$ws.Range("E31").Value = "  -5.30%  "

$ws.Range("D32").Value = "'1.66"
$ws.Range("E32").Value = "  -2.99%  "

$ws.Range("E33").Value = "  +1.89%  "

$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").Value = "'17.62"
$ws.Range("E36").Value = "  -3.60%  "

$ws.Range("E37").Value = "  -5.61%  "

$ws.Range("D38").Value = "'3.90"
$ws.Range("E38").Value = "  -2.35%  "

$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("E40").Value = "  -3.76%  "

$ws.Range("D41").Value = "'0.767"
$ws.Range("E41").Value = "  -4.08%  "

$ws.Range("D42").Value = "'271.23"
$ws.Range("E42").Value = "  -2.75%  "

$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  -5.19%  "

$ws.Range("D44").Value = "'4.86"
$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("D45").Value = "'0.580"
$ws.Range("E45").Value = "  -2.68%  "

$ws.Range("D46").Value = "'0.0905"
$ws.Range("E46").Value = "  -1.07%  "

$ws.Range("D47").Value = "'119.67"
$ws.Range("E47").Value = "  -6.15%  "

$ws.Range("E48").Value = "  -3.23%  "

$ws.Range("D49").Value = "'17.08"
$ws.Range("E49").Value = "  -4.62%  "

$ws.Range("D50").Value = "'0.0209"
$ws.Range("E50").Value = "  -3.63%  "

$ws.Range("D51").Value = "'16.42"
$ws.Range("E51").Value = "  -5.15%  "
